# Continue QA and mex file fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "MAXQDA 12" to "Sheet1"
$ws.Name = "Sheet1"

# 2. Normalize the "Creation date" text values (column M) from
#    MM/DD/YYYY HH:MM:SS to M/D/YY HH:MM:SS for all existing data rows.
$lastRow = 150
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 13)
    $text = $cell.Text
    if ($text -match '^(\d{1,2})/(\d{1,2})/(\d{4}) (\d{2}:\d{2}:\d{2})$') {
        $mm = [int]$matches[1]
        $dd = [int]$matches[2]
        $yyyy = $matches[3]
        $tm = $matches[4]
        $yy = $yyyy.Substring(2,2)
        $cell.Value = "$mm/$dd/$yy $tm"
    }
}

# 3. Append three new coded-segment rows (151-153) for the newly-coded
#    "Location:Hospital name" / "Location:City" / "Location:Country"
#    segments on document 15902 ("Sydney Eye Hospital" / "Sydney" /
#    "Australia"), matching the look/formatting of the existing rows.
$newRows = @(151, 152, 153)
foreach ($r in $newRows) {
    $ws.Range("A150:M150").Copy() | Out-Null
    $ws.Range("A$r`:M$r").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
    $ws.Rows.Item($r).RowHeight = 16
}

# Row 151: Location:Hospital name
$ws.Cells.Item(151, 1).Value = "●"
$ws.Cells.Item(151, 4).NumberFormat = "@"
$ws.Cells.Item(151, 4).Value = "15902"
$ws.Cells.Item(151, 5).Value = "Location:Hospital name"
$ws.Cells.Item(151, 6).Value = "1: 5402"
$ws.Cells.Item(151, 7).Value = "1: 5420"
$ws.Cells.Item(151, 8).Value = 0
$ws.Cells.Item(151, 9).Value = "Sydney Eye Hospital"
$ws.Cells.Item(151, 10).Value = 19
$ws.Cells.Item(151, 11).Value = 0.12952484831958552
$ws.Cells.Item(151, 12).Value = "emmamendelsohn"
$ws.Cells.Item(151, 13).Value = "8/22/19 14:16:18"

# Row 152: Location:City
$ws.Cells.Item(152, 1).Value = "●"
$ws.Cells.Item(152, 4).NumberFormat = "@"
$ws.Cells.Item(152, 4).Value = "15902"
$ws.Cells.Item(152, 5).Value = "Location:City"
$ws.Cells.Item(152, 6).Value = "1: 5423"
$ws.Cells.Item(152, 7).Value = "1: 5428"
$ws.Cells.Item(152, 8).Value = 0
$ws.Cells.Item(152, 9).Value = "Sydney"
$ws.Cells.Item(152, 10).Value = 6
$ws.Cells.Item(152, 11).Value = 0.040902583679869112
$ws.Cells.Item(152, 12).Value = "emmamendelsohn"
$ws.Cells.Item(152, 13).Value = "8/22/19 14:16:22"

# Row 153: Location:Country
$ws.Cells.Item(153, 1).Value = "●"
$ws.Cells.Item(153, 4).NumberFormat = "@"
$ws.Cells.Item(153, 4).Value = "15902"
$ws.Cells.Item(153, 5).Value = "Location:Country"
$ws.Cells.Item(153, 6).Value = "1: 5431"
$ws.Cells.Item(153, 7).Value = "1: 5439"
$ws.Cells.Item(153, 8).Value = 0
$ws.Cells.Item(153, 9).Value = "Australia"
$ws.Cells.Item(153, 10).Value = 9
$ws.Cells.Item(153, 11).Value = 0.061353875519803668
$ws.Cells.Item(153, 12).Value = "emmamendelsohn"
$ws.Cells.Item(153, 13).Value = "8/22/19 14:16:28"

# Restore the "Document name" column's original style (General format,
# shared with the rest of column D) now that the text type has been set.
foreach ($r in $newRows) {
    $ws.Range("D150:D150").Copy() | Out-Null
    $ws.Range("D$r`:D$r").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
    $ws.Rows.Item($r).RowHeight = 16
}
